# Daily attendance processing - reorders the "Recorded By" (column G) values
# so that an exact "System" entry is always listed first, preserving the
# relative order of any remaining entries. When no exact "System" entry is
# present the entries are sorted alphabetically instead.
#
# NOTE: loop/helper variables below are deliberately given unique names
# (prefixed per-scope) because this host does not give functions their own
# variable scope for loop counters - a `for ($i = ...)` inside a helper
# would otherwise clobber an `$i` used by an outer loop.

function Test-ExactEquals($valA, $valB) {
    if ($valA.Length -ne $valB.Length) { return $false }
    $valAChars = $valA.ToCharArray()
    $valBChars = $valB.ToCharArray()
    for ($charIdx = 0; $charIdx -lt $valAChars.Length; $charIdx++) {
        if ([int]$valAChars[$charIdx] -ne [int]$valBChars[$charIdx]) { return $false }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($rowNum = 2; $rowNum -le $lastDataRow; $rowNum++) {
    $gCell = $ws.Cells.Item($rowNum, 7)
    $rawValue = $gCell.Value2

    if ($null -eq $rawValue) { continue }
    if ($rawValue -eq "") { continue }

    $nameParts = @($rawValue -split ", ")
    if ($nameParts.Count -le 1) { continue }

    $systemPos = -1
    for ($partIdx = 0; $partIdx -lt $nameParts.Count; $partIdx++) {
        if (Test-ExactEquals $nameParts[$partIdx] "System") {
            $systemPos = $partIdx
            break
        }
    }

    if ($systemPos -ge 0) {
        $remainingParts = @()
        for ($partIdx2 = 0; $partIdx2 -lt $nameParts.Count; $partIdx2++) {
            if ($partIdx2 -ne $systemPos) { $remainingParts += $nameParts[$partIdx2] }
        }
        $orderedParts = @("System") + $remainingParts
    } else {
        $orderedParts = @($nameParts | Sort-Object)
    }

    $newValue = ($orderedParts -join ", ")

    if ($newValue -ne $rawValue) {
        $gCell.Value = $newValue
    }
}
